# Adds splashscreen logic and loca
#
# The GameName entry ("Kid Game") is replaced with a localized title
# ("Gnome and Owl" / "Gnom und Eule" / "Krasnoludek i Sowa"), and three new
# rows are inserted right below it for the individual title words
# (Gnome/And/Owl) used by the new splashscreen animation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# Insert three fresh rows right after the GameName row (row 12) to hold the
# new per-word localization keys.
$ws.Rows("13:15").Insert() | Out-Null

# GameName's translations change from the placeholder "Kid Game" to the
# real, localized title of the game.
$ws.Range("B12").Value2 = "Gnome and Owl"
$ws.Range("C12").Value2 = "Gnom und Eule"
$ws.Range("D12").Value2 = "Krasnoludek i Sowa"

# New keys/values for the individual title words (written column-by-column
# to match the shared-string allocation order of the authored workbook).
$ws.Range("A13").Value2 = "Gnome"
$ws.Range("A14").Value2 = "And"
$ws.Range("A15").Value2 = "Owl"

$ws.Range("B13").Value2 = "Gnome"
$ws.Range("B14").Value2 = "And"
$ws.Range("B15").Value2 = "Owl"

$ws.Range("C13").Value2 = "Gnom"
$ws.Range("C14").Value2 = "und"
$ws.Range("C15").Value2 = "Eule"

$ws.Range("D13").Value2 = "Krasnoludek"
$ws.Range("D14").Value2 = "i"
$ws.Range("D15").Value2 = "Sowa"

# Leave the selection where the author last left it when saving.
$ws.Range("D16").Select() | Out-Null
